# Update ticket-count (F) figures and sold-out status (G) for the
# 广州-漫展信息 workbook, matching the upstream gh-pages data refresh
# generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 923
$ws.Range("F6").Value = 161
$ws.Range("F7").Value = 946
$ws.Range("F8").Value = 741
$ws.Range("F9").Value = 191
$ws.Range("F12").Value = 773
$ws.Range("F13").Value = 256
$ws.Range("F16").Value = 1296
$ws.Range("F17").Value = 114
$ws.Range("F18").Value = 426
$ws.Range("F19").Value = 1096
$ws.Range("F20").Value = 2795
$ws.Range("F21").Value = 1290
$ws.Range("F22").Value = 651
$ws.Range("F24").Value = 1245
$ws.Range("F26").Value = 970
$ws.Range("F28").Value = 1009
$ws.Range("F29").Value = 19
$ws.Range("F31").Value = 1323

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 512
$ws.Range("G3").Value = "已售罄"
$ws.Range("F4").Value = 349
$ws.Range("F12").Value = 13

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 720

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 720
$ws.Range("F7").Value = 512
$ws.Range("G7").Value = "已售罄"
$ws.Range("F8").Value = 512
$ws.Range("G8").Value = "已售罄"
$ws.Range("F9").Value = 349
$ws.Range("F12").Value = 923
$ws.Range("F13").Value = 161
$ws.Range("F15").Value = 946
$ws.Range("F16").Value = 741
$ws.Range("F17").Value = 191
$ws.Range("F24").Value = 13
$ws.Range("F25").Value = 773
$ws.Range("F26").Value = 256
$ws.Range("F29").Value = 1296
$ws.Range("F30").Value = 114
$ws.Range("F31").Value = 426
$ws.Range("F32").Value = 1096
$ws.Range("F33").Value = 2795
$ws.Range("F34").Value = 1290
$ws.Range("F35").Value = 651
$ws.Range("F37").Value = 1245
$ws.Range("F41").Value = 970
$ws.Range("F43").Value = 1009
$ws.Range("F44").Value = 19
$ws.Range("F46").Value = 1323
